# Revert "Test 2 Laptop": remove the second paragraph entirely
# (including its paragraph mark), leaving only "Test 1 desktop".
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Test 2 Laptop") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
